# Apply the updated cryptocurrency price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.467.57"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "1.849.14"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").Value = "'0.9978"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'240.32"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.6297"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.9988"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "'0.07481"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("E9").Value = "  -0.36%  "
$ws.Range("D10").Value = "'24.62"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "'0.07739"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "'5.014"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'0.6805"
$ws.Range("E14").Value = "  -0.13%  "
$ws.Range("D15").Value = "'0.00001050"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'82.20"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "'6.219"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "29.473.22"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'229.52"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").Value = "'12.38"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'0.9994"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'7.557"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'0.9986"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "'159.19"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'8.541"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'0.1372"
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("D27").Value = "'17.55"
$ws.Range("D28").Value = "'0.06517"
$ws.Range("E28").Value = "  +15.64%  "
$ws.Range("D29").Value = "'1.417"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'1.482"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.102"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.103"
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "'0.6978"
$ws.Range("E35").Value = "  -0.49%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("D38").Value = "1.265.25"
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "'2.840"
$ws.Range("E39").Value = "  +4.21%  "
$ws.Range("D40").Value = "'6.803"
$ws.Range("E40").Value = "  +6.07%  "
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("D42").Value = "'0.9990"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "2.004.77"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "'101.22"
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'66.18"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D47").Value = "'7.099"
$ws.Range("E47").Value = "  -0.91%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1164"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.022"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").Value = "'0.3953"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000114"
$ws.Range("E51").Value = "  -1.40%  "
